$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.413144
$ws.Range("H2").Value = 1.239432
$ws.Range("I2").Value = 0.4553782032534783
$ws.Range("J2").Value = 0.4553782032534782
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 0.8374905372746666
$ws.Range("R2").Value = 7.537414835472
$ws.Range("S2").Value = 0.00300380421354898
$ws.Range("T2").Value = 0.003003804213548979

$ws.Range("G3").Value = 0.413144
$ws.Range("H3").Value = 1.239432
$ws.Range("I3").Value = 0.4553782032534783
$ws.Range("J3").Value = 0.4553782032534782
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 105.9484260826347
$ws.Range("R3").Value = 953.535834743712
$ws.Range("S3").Value = 0.3800022979621162
$ws.Range("T3").Value = 0.3800022979621161

$ws.Range("G4").Value = 0.413144
$ws.Range("H4").Value = 1.239432
$ws.Range("I4").Value = 0.4553782032534783
$ws.Range("J4").Value = 0.4553782032534782
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 20.17806271858934
$ws.Range("R4").Value = 181.602564467304
$ws.Range("S4").Value = 0.07237210107781318
$ws.Range("T4").Value = 0.07237210107781314

$ws.Range("I5").Value = 0.3895918235379703
$ws.Range("J5").Value = 0.3895918235379702
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 0.7165021586046667
$ws.Range("R5").Value = 6.448519427442001
$ws.Range("S5").Value = 0.002569858532416806
$ws.Range("T5").Value = 0.002569858532416805

$ws.Range("I6").Value = 0.3895918235379703
$ws.Range("J6").Value = 0.3895918235379702
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("Q6").Value = 90.64254771881467
$ws.Range("R6").Value = 815.7829294693321
$ws.Range("S6").Value = 0.3251051261434068
$ws.Range("T6").Value = 0.3251051261434066

$ws.Range("I7").Value = 0.3895918235379703
$ws.Range("J7").Value = 0.3895918235379702
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 17.26303146227434
$ws.Range("R7").Value = 155.367283160469
$ws.Range("S7").Value = 0.06191683886214678
$ws.Range("T7").Value = 0.06191683886214675

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1406516666666666
$ws.Range("H8").Value = 0.421955
$ws.Range("I8").Value = 0.1550299732085515
$ws.Range("J8").Value = 0.1550299732085515
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 0.2851171501588888
$ws.Range("R8").Value = 2.56605435143
$ws.Range("S8").Value = 0.00102262181945283
$ws.Range("T8").Value = 0.00102262181945283

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1406516666666666
$ws.Range("H9").Value = 0.421955
$ws.Range("I9").Value = 0.1550299732085515
$ws.Range("J9").Value = 0.1550299732085515
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 36.06931895230888
$ws.Range("R9").Value = 324.62387057078
$ws.Range("S9").Value = 0.1293688315588146
$ws.Range("T9").Value = 0.1293688315588146

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1406516666666666
$ws.Range("H10").Value = 0.421955
$ws.Range("I10").Value = 0.1550299732085515
$ws.Range("J10").Value = 0.1550299732085515
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 6.869464766459444
$ws.Range("R10").Value = 61.825182898135
$ws.Range("S10").Value = 0.02463851983028408
$ws.Range("T10").Value = 0.02463851983028407

Write-Output "updated cells"